$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.370506167411804
$ws.Range("B1").Value = 2.729767799377441
$ws.Range("C1").Value = 3.410092353820801
$ws.Range("D1").Value = 3.412843704223633
$ws.Range("E1").Value = 1.675902247428894
